$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2083333333333333
$ws.Range("C2").Value = 0.5416666666666666
$ws.Range("P2").Value = 0.1666666666666667
$ws.Range("S2").Value = 0.08333333333333333

# Row 3
$ws.Range("P3").Value = 0.7692307692307693
$ws.Range("S3").Value = 0.2307692307692308

# Row 4
$ws.Range("P4").Value = 1

# Row 6
$ws.Range("B6").Value = 0.25
$ws.Range("J6").Value = 0.5
$ws.Range("R6").Value = 0.125
$ws.Range("S6").Value = 0.125

# Row 7
$ws.Range("B7").Value = 0.5
$ws.Range("F7").Value = 0.5

# Row 8
$ws.Range("D8").Value = 0.0625
$ws.Range("F8").Value = 0.125
$ws.Range("J8").Value = 0.1875
$ws.Range("Q8").Value = 0.125
$ws.Range("R8").Value = 0.125
$ws.Range("S8").Value = 0.375

# Row 9
$ws.Range("B9").Value = 0.25
$ws.Range("J9").Value = 0.25
$ws.Range("S9").Value = 0.5

# Row 10
$ws.Range("B10").Value = 0.2222222222222222
$ws.Range("D10").Value = 0.01587301587301587
$ws.Range("F10").Value = 0.01587301587301587
$ws.Range("J10").Value = 0.07936507936507936
$ws.Range("Q10").Value = 0.1746031746031746
$ws.Range("R10").Value = 0.1587301587301587
$ws.Range("S10").Value = 0.3333333333333333

# Row 11
$ws.Range("G11").Value = 0.2
$ws.Range("J11").Value = 0.2
$ws.Range("K11").Value = 0.4
$ws.Range("L11").Value = 0.2

# Row 12
$ws.Range("G12").Value = 1

# Row 15
$ws.Range("I15").Value = 0.1428571428571428
$ws.Range("J15").Value = 0.4285714285714285
$ws.Range("K15").Value = 0.2857142857142857
$ws.Range("S15").Value = 0.1428571428571428

# Row 16
$ws.Range("I16").Value = 0.1333333333333333
$ws.Range("J16").Value = 0.7333333333333333
$ws.Range("O16").Value = 0.1333333333333333

# Row 17
$ws.Range("F17").Value = 0.07692307692307693
$ws.Range("H17").Value = 0.1538461538461539
$ws.Range("J17").Value = 0.6153846153846154
$ws.Range("S17").Value = 0.1538461538461539

# Row 18
$ws.Range("H18").Value = 0.2307692307692308
$ws.Range("J18").Value = 0.5384615384615384
$ws.Range("O18").Value = 0.07692307692307693
$ws.Range("S18").Value = 0.1538461538461539

# Row 19
$ws.Range("F19").Value = 0.0425531914893617
$ws.Range("H19").Value = 0.2340425531914894
$ws.Range("I19").Value = 0.02127659574468085
$ws.Range("J19").Value = 0.425531914893617
$ws.Range("K19").Value = 0.02127659574468085
$ws.Range("M19").Value = 0.02127659574468085
$ws.Range("O19").Value = 0.0851063829787234
$ws.Range("S19").Value = 0.148936170212766

